$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '61.713.23'
$ws.Range("E2").Value = '  -2.43%  '
$ws.Range("D3").Value = '3.395.99'
$ws.Range("E3").Value = '  -2.34%  '
$ws.Range("E4").Value = '  -0.14%  '
Set-TextValue "D5" '406.95'
$ws.Range("E5").Value = '  -2.22%  '
Set-TextValue "D6" '134.93'
$ws.Range("E6").Value = '  +8.51%  '
Set-TextValue "D7" '0.594'
$ws.Range("E7").Value = '  -1.12%  '
$ws.Range("E8").Value = '  -0.03%  '
Set-TextValue "D9" '0.674'
$ws.Range("E9").Value = '  -0.28%  '
Set-TextValue "D10" '0.122'
$ws.Range("E10").Value = '  -6.94%  '
Set-TextValue "D11" '42.92'
$ws.Range("E11").Value = '  +3.20%  '
$ws.Range("E12").Value = '  -1.12%  '
$ws.Range("D13").Value = '3.910.80'
$ws.Range("E13").Value = '  -2.87%  '
Set-TextValue "D14" '8.44'
$ws.Range("E14").Value = '  -1.74%  '
Set-TextValue "D15" '19.82'
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").Value = '3.396.39'
$ws.Range("E16").Value = '  -2.71%  '
$ws.Range("D17").Value = '61.579.07'
$ws.Range("E17").Value = '  -2.47%  '
$ws.Range("E18").Value = '  -1.22%  '
Set-TextValue "D19" '11.01'
$ws.Range("E19").Value = '  -1.14%  '
Set-TextValue "D20" '0.0000128'
$ws.Range("E20").Value = '  -6.08%  '
$ws.Range("E21").Value = '  -4.04%  '
Set-TextValue "D22" '85.07'
$ws.Range("E22").Value = '  +2.41%  '
Set-TextValue "D23" '315.50'
$ws.Range("E23").Value = '  -0.58%  '
Set-TextValue "D24" '12.85'
$ws.Range("E24").Value = '  -0.93%  '
Set-TextValue "D25" '3.15'
$ws.Range("E25").Value = '  -1.45%  '
Set-TextValue "D26" '4.79'
$ws.Range("E26").Value = '  +11.01%  '
Set-TextValue "D27" '8.37'
$ws.Range("E27").Value = '  +5.87%  '
Set-TextValue "D28" '29.64'
$ws.Range("E28").Value = '  -4.67%  '
Set-TextValue "D29" '7.60'
$ws.Range("E29").Value = '  -2.88%  '
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("E31").Value = '  +3.61%  '
Set-TextValue "D32" '0.171'
$ws.Range("E32").Value = '  -2.62%  '
Set-TextValue "D33" '11.42'
$ws.Range("E33").Value = '  -1.95%  '
Set-TextValue "D34" '1.00'
$ws.Range("E34").Value = '  -0.51%  '
Set-TextValue "D35" '41.44'
$ws.Range("E35").Value = '  -1.51%  '
Set-TextValue "D36" '0.0484'
$ws.Range("E36").Value = '  -1.40%  '
Set-TextValue "D37" '51.66'
$ws.Range("E37").Value = '  -0.95%  '
Set-TextValue "D38" '0.996'
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("E39").Value = '  -1.88%  '
Set-TextValue "D40" '2.95'
$ws.Range("E40").Value = '  -3.32%  '
Set-TextValue "D41" '140.36'
$ws.Range("E41").Value = '  +3.85%  '
Set-TextValue "D42" '1.99'
$ws.Range("E42").Value = '  -1.36%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue "D43" '0.298'
$ws.Range("E43").Value = '  +4.50%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D44" '0.124'
$ws.Range("E44").Value = '  -1.87%  '
Set-TextValue "D45" '4.03'
$ws.Range("E45").Value = '  +2.69%  '
$ws.Range("E46").Value = '  -3.52%  '
$ws.Range("E47").Value = '  -0.89%  '
Set-TextValue "D48" '21.39'
$ws.Range("E48").Value = '  -3.72%  '
$ws.Range("D49").Value = '2.118.96'
$ws.Range("E49").Value = '  -3.11%  '
$ws.Range("E50").Value = '  -4.92%  '
Set-TextValue "D51" '1.91'
$ws.Range("E51").Value = '  -0.24%  '
